$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = "12 x 3 =|20 - 9 =|56 + 9 =|40 ÷ 10 ="
$ws.Range("C13").Value = "5 ÷ 2 |8 ÷ 8|11 ÷ 3|19 ÷ 4"
$ws.Range("C16").Value = "10 + 5 =|6 - 3 =|9 x 3 =!9 ÷ 3 |25 - 10|15 + 12"
$ws.Range("C17").Value = "3#3 = 9|1#2 = 3|10#2 = 5|7#3 = 4!+|-|x|÷"

$ws.Range("C18").Select()
